$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This shared string is referenced by the Status column (B) and, on the
#    Overview sheet, also by column C, for the two tracked files (rows 2-3).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B2").Value = $newStatus
$wsZhCn.Range("B3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B2").Value = $newStatus
$wsDeDe.Range("B3").Value = $newStatus

# Hyperlink-style font used for the file-name columns (matches the "HyperLink"
# cell style already used by columns A and C: underlined, cornflowerblue).
$linkColor = 15570276

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: populate the "Latest Target File" (E) and
#    "Latest Handback File" (F) columns for both tracked rows, with
#    hyperlinks that mirror the existing Source/Handoff links of the row.
# ---------------------------------------------------------------------------
$zhE2Url = "https://github.com/OpenLocalizationTest/oltest/blob/78abb8eea48daaa7edfcc9b97946f27fb85a00ee/e2e/84b09259-7555-4c1e-b1de-2f97f75eef95.md"
$wsZhCn.Range("E2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), $zhE2Url, "", "", "84b09259-7555-4c1e-b1de-2f97f75eef95.md")
$wsZhCn.Range("E2").Font.Underline = $true
$wsZhCn.Range("E2").Font.Color = $linkColor

$zhF2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a5944192b7a351018b70b9a3eef671b1d39445c5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/84b09259-7555-4c1e-b1de-2f97f75eef95.ebc2fbcc84a15bf51808aaae1467fa9b88a4ceb3.zh-cn.xlf"
$wsZhCn.Range("F2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.ebc2fbcc84a15bf51808aaae1467fa9b88a4ceb3.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhF2Url, "", "", "84b09259-7555-4c1e-b1de-2f97f75eef95.ebc2fbcc84a15bf51808aaae1467fa9b88a4ceb3.zh-cn.xlf")
$wsZhCn.Range("F2").Font.Underline = $true
$wsZhCn.Range("F2").Font.Color = $linkColor

$zhE3Url = "https://github.com/OpenLocalizationTest/oltest/blob/78abb8eea48daaa7edfcc9b97946f27fb85a00ee/e2e/a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.md"
$wsZhCn.Range("E3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), $zhE3Url, "", "", "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.md")
$wsZhCn.Range("E3").Font.Underline = $true
$wsZhCn.Range("E3").Font.Color = $linkColor

$zhF3Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a5944192b7a351018b70b9a3eef671b1d39445c5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.fe500b8130d0f1a52ea71db28b4502d3da31c4d8.zh-cn.xlf"
$wsZhCn.Range("F3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.fe500b8130d0f1a52ea71db28b4502d3da31c4d8.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhF3Url, "", "", "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.fe500b8130d0f1a52ea71db28b4502d3da31c4d8.zh-cn.xlf")
$wsZhCn.Range("F3").Font.Underline = $true
$wsZhCn.Range("F3").Font.Color = $linkColor

# ---------------------------------------------------------------------------
# 3. de-de sheet: same as above, plus the handback datetime (G) is now known
#    for both rows (previously the "0001-01-01 00:00:00" placeholder).
# ---------------------------------------------------------------------------
$deE2Url = "https://github.com/OpenLocalizationTest/oltest/blob/78abb8eea48daaa7edfcc9b97946f27fb85a00ee/e2e/84b09259-7555-4c1e-b1de-2f97f75eef95.md"
$wsDeDe.Range("E2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), $deE2Url, "", "", "84b09259-7555-4c1e-b1de-2f97f75eef95.md")
$wsDeDe.Range("E2").Font.Underline = $true
$wsDeDe.Range("E2").Font.Color = $linkColor

$deF2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/be62f80592bcbdb48b5637e3aced34cd03e6daec/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/84b09259-7555-4c1e-b1de-2f97f75eef95.ebc2fbcc84a15bf51808aaae1467fa9b88a4ceb3.de-de.xlf"
$wsDeDe.Range("F2").Value = "84b09259-7555-4c1e-b1de-2f97f75eef95.ebc2fbcc84a15bf51808aaae1467fa9b88a4ceb3.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deF2Url, "", "", "84b09259-7555-4c1e-b1de-2f97f75eef95.ebc2fbcc84a15bf51808aaae1467fa9b88a4ceb3.de-de.xlf")
$wsDeDe.Range("F2").Font.Underline = $true
$wsDeDe.Range("F2").Font.Color = $linkColor

$wsDeDe.Range("G2").Value = "2016-01-26 09:06:15"

$deE3Url = "https://github.com/OpenLocalizationTest/oltest/blob/78abb8eea48daaa7edfcc9b97946f27fb85a00ee/e2e/a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.md"
$wsDeDe.Range("E3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), $deE3Url, "", "", "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.md")
$wsDeDe.Range("E3").Font.Underline = $true
$wsDeDe.Range("E3").Font.Color = $linkColor

$deF3Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/be62f80592bcbdb48b5637e3aced34cd03e6daec/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.fe500b8130d0f1a52ea71db28b4502d3da31c4d8.de-de.xlf"
$wsDeDe.Range("F3").Value = "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.fe500b8130d0f1a52ea71db28b4502d3da31c4d8.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deF3Url, "", "", "a7ba330e-2f4c-4ebd-b97b-562e88ee6c29.fe500b8130d0f1a52ea71db28b4502d3da31c4d8.de-de.xlf")
$wsDeDe.Range("F3").Font.Underline = $true
$wsDeDe.Range("F3").Font.Color = $linkColor

$wsDeDe.Range("G3").Value = "2016-01-26 09:06:15"
